$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 253.26666
$ws.Range("I33").Value = 264.2143
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 264.2143
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = -35.21429999999998
$ws.Range("N33").Value = -558
$ws.Range("H51").Value = 10591
$ws.Range("J51").Value = 10591
$ws.Range("L51").Value = 10591
$ws.Range("N51").Value = -11559
$ws.Range("H53").Value = 470.52
$ws.Range("J53").Value = 696.8333
$ws.Range("L53").Value = 696.8333
$ws.Range("N53").Value = -1970.8333
$ws.Range("H62").Value = 3109.6
$ws.Range("J62").Value = 5250
$ws.Range("L62").Value = 5250
$ws.Range("N62").Value = -6498
$ws.Range("H65").Value = 3109.6
$ws.Range("J65").Value = 5250
$ws.Range("L65").Value = 26250
$ws.Range("N65").Value = -32490
$ws.Range("H93").Value = 29881.312
$ws.Range("J93").Value = 29881.312
$ws.Range("L93").Value = 29881.312
$ws.Range("N93").Value = -34873.31200000001
$ws.Range("H116").Value = 918900.5600000001
$ws.Range("I116").Value = 3336666.8
$ws.Range("K116").Value = 3336666.8
$ws.Range("M116").Value = -3333224.8
$ws.Range("H129").Value = 857.0599999999999
$ws.Range("J129").Value = 871.91754
$ws.Range("L129").Value = 2615.75262
$ws.Range("N129").Value = -12615.75262
$ws.Range("H138").Value = 2723.5
$ws.Range("I138").Value = 1195.3334
$ws.Range("J138").Value = 2874.6375
$ws.Range("K138").Value = 3586.0002
$ws.Range("L138").Value = 8623.912499999999
$ws.Range("M138").Value = 1553.9998
$ws.Range("N138").Value = -18903.9125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 405
$ws.Range("I2").Value = 405
$ws.Range("K2").Value = 405
$ws.Range("M2").Value = -292
$ws.Range("H24").Value = 22230.6
$ws.Range("J24").Value = 22230.6
$ws.Range("L24").Value = 22230.6
$ws.Range("N24").Value = -22978.6
$ws.Range("H46").Value = 5000
$ws.Range("I46").Value = 5000
$ws.Range("K46").Value = 5000
$ws.Range("M46").Value = -4681
$ws.Range("H100").Value = 22230.6
$ws.Range("J100").Value = 22230.6
$ws.Range("L100").Value = 22230.6
$ws.Range("N100").Value = -24394.6
$ws.Range("H116").Value = 405
$ws.Range("I116").Value = 405
$ws.Range("K116").Value = 405
$ws.Range("M116").Value = 1889
$ws.Range("H122").Value = 6240.6665
$ws.Range("I122").Value = 1861
$ws.Range("K122").Value = 5583
$ws.Range("M122").Value = -3133
$ws.Range("H131").Value = 40514
$ws.Range("J131").Value = 40514
$ws.Range("L131").Value = 40514
$ws.Range("N131").Value = -50594
$ws.Range("H133").Value = 19533.143
$ws.Range("J133").Value = 19533.143
$ws.Range("L133").Value = 19533.143
$ws.Range("N133").Value = -24593.143
$ws.Range("H139").Value = 39805.445
$ws.Range("J139").Value = 39805.445
$ws.Range("L139").Value = 39805.445
$ws.Range("N139").Value = -50085.445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 405
$ws.Range("I3").Value = 405
$ws.Range("K3").Value = 405
$ws.Range("M3").Value = -291
$ws.Range("H95").Value = 35050
$ws.Range("J95").Value = 35050
$ws.Range("L95").Value = 35050
$ws.Range("N95").Value = -40542
$ws.Range("H134").Value = 2482.625
$ws.Range("I134").Value = 1898.56
$ws.Range("K134").Value = 5695.68
$ws.Range("M134").Value = -3160.68
$ws.Range("H138").Value = 40681.035
$ws.Range("J138").Value = 40681.035
$ws.Range("L138").Value = 40681.035
$ws.Range("N138").Value = -50961.035

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 21457.5
$ws.Range("I55").Value = 10500
$ws.Range("J55").Value = 28762.5
$ws.Range("K55").Value = 10500
$ws.Range("L55").Value = 28762.5
$ws.Range("M55").Value = -10185
$ws.Range("N55").Value = -29392.5
$ws.Range("H58").Value = 3405.1428
$ws.Range("I58").Value = 1947.4546
$ws.Range("J58").Value = 8750
$ws.Range("K58").Value = 1947.4546
$ws.Range("L58").Value = 8750
$ws.Range("M58").Value = -1744.4546
$ws.Range("N58").Value = -9156
$ws.Range("H68").Value = 84999.25
$ws.Range("J68").Value = 84999.25
$ws.Range("L68").Value = 84999.25
$ws.Range("N68").Value = -86497.25
$ws.Range("H71").Value = 84999.25
$ws.Range("J71").Value = 84999.25
$ws.Range("L71").Value = 254997.75
$ws.Range("N71").Value = -262485.75
$ws.Range("H132").Value = 2255.1353
$ws.Range("I132").Value = 1452.4615
$ws.Range("J132").Value = 4152.364
$ws.Range("K132").Value = 4357.3845
$ws.Range("L132").Value = 12457.092
$ws.Range("M132").Value = -1827.3845
$ws.Range("N132").Value = -17517.092
$ws.Range("H136").Value = 3405.1428
$ws.Range("I136").Value = 1947.4546
$ws.Range("J136").Value = 8750
$ws.Range("K136").Value = 5842.3638
$ws.Range("L136").Value = 26250
$ws.Range("M136").Value = -3292.3638
$ws.Range("N136").Value = -31350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1406.7693
$ws.Range("I41").Value = 379.6
$ws.Range("J41").Value = 2048.75
$ws.Range("K41").Value = 1138.8
$ws.Range("L41").Value = 6146.25
$ws.Range("M41").Value = -800.8000000000002
$ws.Range("N41").Value = -6822.25
$ws.Range("H131").Value = 706.2371000000001
$ws.Range("J131").Value = 796.525
$ws.Range("L131").Value = 2389.575
$ws.Range("N131").Value = -12469.575

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 34900
$ws.Range("J105").Value = 34900
$ws.Range("L105").Value = 34900
$ws.Range("N105").Value = -41888
$ws.Range("H122").Value = 12400
$ws.Range("J122").Value = 14250
$ws.Range("L122").Value = 42750
$ws.Range("N122").Value = -47650
$ws.Range("H124").Value = 41824
$ws.Range("J124").Value = 41824
$ws.Range("L124").Value = 41824
$ws.Range("N124").Value = -51644
$ws.Range("H126").Value = 3313.42
$ws.Range("I126").Value = 2932.8933
$ws.Range("J126").Value = 4455
$ws.Range("K126").Value = 8798.679900000001
$ws.Range("L126").Value = 13365
$ws.Range("M126").Value = -6328.679900000001
$ws.Range("N126").Value = -18305
$ws.Range("H130").Value = 48958.75
$ws.Range("J130").Value = 48958.75
$ws.Range("L130").Value = 48958.75
$ws.Range("N130").Value = -58998.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2715.111
$ws.Range("I61").Value = 2822.8333
$ws.Range("J61").Value = 2499.6667
$ws.Range("K61").Value = 2822.8333
$ws.Range("L61").Value = 2499.6667
$ws.Range("M61").Value = -2620.8333
$ws.Range("N61").Value = -2903.6667
$ws.Range("H104").Value = 24999.666
$ws.Range("J104").Value = 24999.666
$ws.Range("L104").Value = 24999.666
$ws.Range("N104").Value = -31987.666
$ws.Range("H113").Value = 2715.111
$ws.Range("I113").Value = 2822.8333
$ws.Range("J113").Value = 2499.6667
$ws.Range("K113").Value = 2822.8333
$ws.Range("L113").Value = 2499.6667
$ws.Range("M113").Value = -652.8332999999998
$ws.Range("N113").Value = -6839.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 18000
$ws.Range("J14").Value = 18000
$ws.Range("L14").Value = 18000
$ws.Range("N14").Value = -18336
$ws.Range("H108").Value = 34642
$ws.Range("J108").Value = 34642
$ws.Range("L108").Value = 34642
$ws.Range("N108").Value = -42322
$ws.Range("H132").Value = 8132043.5
$ws.Range("I132").Value = 1408.0667
$ws.Range("J132").Value = 30306504
$ws.Range("K132").Value = 4224.2001
$ws.Range("L132").Value = 90919512
$ws.Range("M132").Value = -1694.2001
$ws.Range("N132").Value = -90924572
$ws.Range("H136").Value = 4622.1035
$ws.Range("I136").Value = 2774.4614
$ws.Range("J136").Value = 6123.3125
$ws.Range("K136").Value = 8323.3842
$ws.Range("L136").Value = 18369.9375
$ws.Range("M136").Value = -6123.3842
$ws.Range("N136").Value = -23469.9375
